$d = $word.ActiveDocument

$replacements = @(
    @{Old = "91×97=8827"; New = "14×70=980"},
    @{Old = "86×43=3698"; New = "59×28=1652"},
    @{Old = "22×47=1034"; New = "68×81=5508"},
    @{Old = "65×82=5330"; New = "71×36=2556"},
    @{Old = "82×15=1230"; New = "58×80=4640"},
    @{Old = "21×84=1764"; New = "79×54=4266"},
    @{Old = "80×23=1840"; New = "43×85=3655"},
    @{Old = "56×44=2464"; New = "33×36=1188"},
    @{Old = "92×98=9016"; New = "53×45=2385"},
    @{Old = "25×36=900"; New = "42×67=2814"},
    @{Old = "98×41=4018"; New = "47×17=799"},
    @{Old = "59×80=4720"; New = "82×85=6970"},
    @{Old = "19×35=665"; New = "45×33=1485"},
    @{Old = "52×94=4888"; New = "20×71=1420"},
    @{Old = "88×65=5720"; New = "23×78=1794"},
    @{Old = "67×91=6097"; New = "42×34=1428"},
    @{Old = "15×23=345"; New = "95×66=6270"},
    @{Old = "43×98=4214"; New = "48×50=2400"},
    @{Old = "94×45=4230"; New = "76×11=836"},
    @{Old = "94×94=8836"; New = "19×84=1596"},
    @{Old = "35×65=2275"; New = "59×54=3186"},
    @{Old = "36×30=1080"; New = "68×83=5644"},
    @{Old = "62×17=1054"; New = "68×65=4420"},
    @{Old = "83×22=1826"; New = "36×91=3276"},
    @{Old = "88×78=6864"; New = "29×82=2378"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

Write-Host "Done applying replacements"
